# 989: changes to make new team in ND02 file
#
# The ND02 extract currently contains data for the "London"/"ND01" team.
# This creates a brand-new "II" team (region "London II", provider code
# "ND02", team "WMT (ND02)", officers "Swann II" / "Wright II") and points
# the two data sheets (WMT_Extract, CMS) at it. All other sheets (which
# reference generic lookup codes, not these team-specific strings) are left
# untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "WMT_Extract": Region/Provider/Team/Surname columns -----------
$ws1 = $wb.Worksheets.Item("WMT_Extract")

# Row 2 (Tom Swann's row)
$ws1.Range("C2").Value = "ND02"
$ws1.Range("H2").Value = "Swann II"

# Row 3 (Andy Wright's row) - surname first, to match new-team roster order
$ws1.Range("H3").Value = "Wright II"

$ws1.Range("G2").Value = "WMT (ND02)"
$ws1.Range("B2").Value = "London II"

$ws1.Range("C3").Value = "ND02"
$ws1.Range("G3").Value = "WMT (ND02)"
$ws1.Range("B3").Value = "London II"

$ws1.Activate()
$ws1.Range("E5").Select()

# --- Sheet "CMS": Contact/OM team + provider + staff keys ----------------
$ws2 = $wb.Worksheets.Item("CMS")

# Row 2: contact by Tom Swann (key 1001->1004), about Andy Wright (1002->1005)
$ws2.Range("F2").Value = 1004
$ws2.Range("H2").Value = "WMT (ND02)"
$ws2.Range("I2").Value = "ND02"
$ws2.Range("K2").Value = 1005
$ws2.Range("M2").Value = "WMT (ND02)"
$ws2.Range("N2").Value = "ND02"

# Row 3: contact by Andy Wright (key 1002->1005), about Tom Swann (1001->1004)
$ws2.Range("F3").Value = 1005
$ws2.Range("H3").Value = "WMT (ND02)"
$ws2.Range("I3").Value = "ND02"
$ws2.Range("K3").Value = 1004
$ws2.Range("M3").Value = "WMT (ND02)"
$ws2.Range("N3").Value = "ND02"

$ws2.Activate()
$ws2.Range("I3").Select()
